$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 450
$ws.Range("B3").Value = 200
$ws.Range("B4").Value = 700
$ws.Range("B5").Value = 78
$ws.Range("B6").Value = 100
$ws.Range("B7").Value = 224
$ws.Range("B8").Value = 98
